$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-11-13"

# Update the column header label in I1 ("2022 (through 11-12)" -> "2022 (through 11-13)")
$ws.Range("I1").Value = "2022 (through 11-13)"

# Update the November row's 2022 value (I12: 35 -> 42)
$ws.Range("I12").Value = 42

# Update the Total row's 2022 value (I14: 1433 -> 1440)
$ws.Range("I14").Value = 1440
